$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") for rows 2-6: change date serial 45221 -> 45224
# (2023-10-22 -> 2023-10-25), preserving existing cell formatting.
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
